# Update "Status" (column B) for a batch of games to "Platinado" (Platinum
# trophy achieved), and Undertale's status from "Zerar" to "Zerado".
#
# Rows (1-based, as in the sheet) that change to "Platinado":
#   11 (Days Gone), 13 (Elden Ring), 16 (Ghost of Tsushima), 20 (God of War 3),
#   21 (God of War 4), 22 (God of War 4), 45 (Resident Evil 2),
#   46 (Resident Evil 2), 47 (Resident Evil 3), 48 (Resident Evil 4),
#   49 (Resident Evil 4 (2005)), 50 (Resident Evil 4 (2005)),
#   51 (Resident Evil 7), 53 (Resident Evil Village), 73 (Hollow Knight),
#   75 (Mortal Kombat 1)
#
# Row that changes to "Zerado":
#   71 (Undertale)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$platinadoRows = @(11, 13, 16, 20, 21, 22, 45, 46, 47, 48, 49, 50, 51, 53, 73, 75)
foreach ($r in $platinadoRows) {
    $ws.Cells.Item($r, 2).Value = "Platinado"
}

$ws.Cells.Item(71, 2).Value = "Zerado"

# Scroll the window so the top-left visible cell is A43 (matches saved view state).
$excel.ActiveWindow.ScrollRow = 43
